$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'90.757.24"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "'3.165.55"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'214.46"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "'630.08"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "'0.405"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").Value = "'0.726"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'3.165.79"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").Value = "'0.562"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'0.182"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'90.578.71"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "'5.30"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "'3.762.71"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "'32.31"
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").Value = "'3.187.56"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'3.31"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").Value = "'0.0000213"
$ws.Range("E20").Value = "  +31.87%  "
$ws.Range("D21").Value = "'13.34"
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("D22").Value = "'431.10"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").Value = "'8.41"
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("D24").Value = "'4.96"
$ws.Range("E24").Value = "  -4.59%  "
$ws.Range("D25").Value = "'5.26"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").Value = "'11.61"
$ws.Range("E26").Value = "  -7.07%  "
$ws.Range("D27").Value = "'80.67"
$ws.Range("E27").Value = "  +7.61%  "
$ws.Range("D28").Value = "'3.352.23"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'0.158"
$ws.Range("E30").Value = "  -9.10%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'4.01"
$ws.Range("E32").Value = "  +25.89%  "
$ws.Range("D33").Value = "'8.30"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").Value = "'509.53"
$ws.Range("E34").Value = "  -9.51%  "
$ws.Range("D35").Value = "'6.90"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("D36").Value = "'1.87"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").Value = "'22.22"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'22.34"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "'0.126"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'1.91"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "'0.371"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'147.18"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Value = "'43.89"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "'168.52"
$ws.Range("E47").Value = "  -6.18%  "
$ws.Range("D48").Value = "'0.125"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'0.734"
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("D50").Value = "'24.43"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "'1.20"
$ws.Range("E51").Value = "  -5.04%  "
